# Replace the single empty paragraph with five new paragraphs of text.
# Each paragraph is tagged as English (US) (w:lang="en-US") on both the
# paragraph mark and the run, matching what Word stamps on typed text.
# Three of the five words are "misspelled" or a find spelled word, so Word's
# proofer wraps them in spellStart/spellEnd proofErr markers; "Sathish"
# and "fedrstae" are left unmarked.

$d = $word.ActiveDocument

function New-WordRun([string]$text) {
    $lang = '<w:rPr><w:lang w:val="en-US"/></w:rPr>'
    return '<w:r>' + $lang + '<w:t>' + $text + '</w:t></w:r>'
}

function New-WordParagraph([string]$text, [bool]$flagSpelling) {
    $pPr = '<w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr>'
    $run = New-WordRun $text
    if ($flagSpelling) {
        $runContent = '<w:proofErr w:type="spellStart"/>' + $run + '<w:proofErr w:type="spellEnd"/>'
    } else {
        $runContent = $run
    }
    return '<w:p>' + $pPr + $runContent + '</w:p>'
}

$paragraphs =
    (New-WordParagraph "Sathish"     $false) +
    (New-WordParagraph "Stahish"     $true)  +
    (New-WordParagraph "Sdefsgjssl"  $true)  +
    (New-WordParagraph "Fedrstansls" $true)  +
    (New-WordParagraph "fedrstae"    $false)

$bodyXml = '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + $paragraphs + '</w:body></w:document>'

$flatOpc = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
    '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
    '<pkg:xmlData>' + $bodyXml + '</pkg:xmlData>' +
    '</pkg:part>' +
    '</pkg:package>'

# Target the whole (currently empty) body content -- Start to End -- so the
# lone pre-existing empty paragraph is replaced rather than left behind.
$target = $d.Range($d.Content.Start, $d.Content.End)
$target.InsertXML($flatOpc)
